$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: fill in the transponder board's 9-pin terminal connector row ---

# B16: part number label, formatted like the other "part number" cells in column B
# (e.g. B4, B7, B9, B10) which use a wrapped, small Arial font.
$ws.Range("B4").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B16").Value = "VI0921550000G"

# C16: replace the old placeholder description with the real part description
$ws.Range("C16").Value = "9 pin terminal"

# D16 (quantity) stays 1 - already set

# E16: price per unit
$ws.Range("E16").Value = 3.1

# F16: source URL - match formatting used by the rest of column F (style index 5)
$ws.Range("F16").ReadingOrder = 1
$ws.Range("F16").Value = "https://www.digikey.ca/product-detail/en/amphenol-anytek/VI0921550000G/609-3943-ND/2261377"

# --- Update the saved selection to reflect where the user finished editing ---
$ws.Range("B16").Select() | Out-Null

# --- Page setup: the sheet was set to print in portrait orientation ---
$ws.PageSetup.Orientation = 1
